# Poland IV Liga - base update (29-03-2024)
#
# The source sheet had several match rows whose data (everything except the
# running "id" in column A) needs to be rotated among a small cluster of
# rows that share the same kickoff date, and the very last row (a fixture
# that hadn't been played yet) needs to be removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a cycle of row numbers. For a cycle (r0, r1, ..., rn-1)
# the NEW content of r_i (columns B..AC) becomes the OLD content of r_(i+1),
# wrapping around - i.e. the rows rotate their data one step "up" the list.
$cycles = @(
    @(30, 31),
    @(44, 45),
    @(46, 47),
    @(59, 60, 61),
    @(65, 66),
    @(73, 74),
    @(75, 76),
    @(86, 87),
    @(108, 109, 110),
    @(124, 125)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot every row's current B:AC content before writing anything,
    # since rows within a cycle feed each other.
    $snapshots = @()
    for ($i = 0; $i -lt $n; $i++) {
        $r = $cycle[$i]
        $snapshots += , ($ws.Range("B$r`:AC$r").Value())
    }

    for ($i = 0; $i -lt $n; $i++) {
        $r = $cycle[$i]
        $srcValues = $snapshots[($i + 1) % $n]
        $ws.Range("B$r`:AC$r").Value = $srcValues
    }
}

# The final row (id 127, match 8022647) is removed outright - it hasn't been
# played yet and drops out of this export entirely.
$ws.Rows.Item(129).Delete()
